$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at K (pushes old K,L -> M,N)
$ws.Range("K1:L1").EntireColumn.Insert()

# New header cells for the inserted columns
$ws.Range("K1").Value = "fi"
$ws.Range("K1").Style = "Normal"
$ws.Range("L1").Value = "se"
$ws.Range("L1").Style = "Normal"

# Fill in new translation content
$ws.Range("H2").Value = "test de_DE"
$ws.Range("H2").Style = "Normal"

$ws.Range("I3").Value = "παράδειγμα"
$ws.Range("J3").Value = "2023-11-09"

$ws.Range("I4").Value = "ψάρι"
$ws.Range("J4").Value = "2023-11-09"
